$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (scrape refresh performed by the
# "Updated cryptos list" GitHub Actions workflow on Sat May  4 09:14:39 UTC 2024).
# Columns:
#   B = Coin name, C = coinranking.com link, D = Price, E = Volume(1h) change.
# Price values such as "586.02" or "0.999" must stay as text (the sheet always
# stores Price as text, e.g. "1.00"), so NumberFormat is forced to "@" (Text)
# before assigning any value that Excel would otherwise auto-convert to a number.

$ws.Range("D2").Value = "63.330.17"
$ws.Range("E2").Value = "  +6.66%  "
$ws.Range("D3").Value = "3.117.01"
$ws.Range("E3").Value = "  +4.60%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.02"
$ws.Range("E5").Value = "  +3.60%  "
$ws.Range("E6").Value = "  +5.18%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.109.00"
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +15.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("E11").Value = "  +7.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").Value = "  +3.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +8.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.68"
$ws.Range("E14").Value = "  +6.02%  "
$ws.Range("D16").Value = "3.631.86"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.18"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "63.235.99"
$ws.Range("E18").Value = "  +6.50%  "
$ws.Range("D19").Value = "3.112.46"
$ws.Range("E19").Value = "  +4.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.54"
$ws.Range("E20").Value = "  +7.17%  "
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.728"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +7.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.35"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.17"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.53"
$ws.Range("E27").Value = "  +10.27%  "
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("E29").Value = "  +5.27%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.92"
$ws.Range("E31").Value = "  +11.68%  "
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("E33").Value = "  +4.82%  "
$ws.Range("D34").Value = "0.0₃0869"
$ws.Range("E34").Value = "  +13.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +17.16%  "
$ws.Range("E36").Value = "  +6.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.34"
$ws.Range("E37").Value = "  +21.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.09"
$ws.Range("E38").Value = "  +3.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.55"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "439.55"
$ws.Range("E40").Value = "  +9.59%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").Value = "2.917.97"
$ws.Range("E42").Value = "  +6.63%  "
$ws.Range("E43").Value = "  +5.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.280"
$ws.Range("E44").Value = "  +11.93%  "
$ws.Range("E45").Value = "  +5.53%  "
$ws.Range("E46").Value = "  +8.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.32"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.69"
$ws.Range("E51").Value = "  +6.28%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.97"
$ws.Range("E48").Value = "  -0.10%  "

Write-Output "Applied crypto data update."
